# Add the missing "Eliminar las conjuntas..." bullet right after the
# "Borrar los comentarios..." bullet, inheriting that bullet's list
# formatting (style "Prrafodelista", numId 6, same spacing).

$d = $word.ActiveDocument

$anchorText = "Borrar los comentarios que considere inadecuados."
$newText = "Eliminar las conjuntas que considere inadecuadas"

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq $anchorText) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find anchor paragraph ending in: $anchorText"
}

# InsertParagraphAfter() splits in a new paragraph immediately following
# $target, inheriting $target's paragraph formatting (style, numbering,
# spacing, etc.) automatically - exactly what we want for a new list item.
$target.Range.InsertParagraphAfter() | Out-Null

$newPara = $target.Next()
$newPara.Range.Text = $newText
